$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This update refreshes the scraped crypto price/volume table.
# Column D ("Price") values are plain decimal-looking strings (e.g. "1.00",
# "35.20") that must stay as literal text -- otherwise Excel would silently
# coerce them to numbers and drop meaningful trailing zeros. Force those
# cells to Text format before assigning. Columns B/C/E (name, link, percent
# strings) are never numeric-looking, so a plain value assignment is enough.

$priceUpdates = @{
    "D2" = "47.404.01"
    "D3" = "2.486.61"
    "D5" = "322.65"
    "D6" = "105.62"
    "D8" = "1.00"
    "D9" = "0.543"
    "D10" = "38.09"
    "D11" = "0.0812"
    "D13" = "18.39"
    "D14" = "7.18"
    "D15" = "2.880.23"
    "D16" = "2.492.36"
    "D17" = "0.846"
    "D18" = "47.320.95"
    "D19" = "12.77"
    "D20" = "6.56"
    "D21" = "0.0₃0937"
    "D22" = "70.62"
    "D24" = "251.47"
    "D26" = "26.15"
    "D28" = "2.26"
    "D29" = "10.02"
    "D30" = "35.20"
    "D32" = "49.44"
    "D33" = "19.63"
    "D35" = "0.0787"
    "D37" = "1.97"
    "D38" = "4.62"
    "D39" = "2.98"
    "D41" = "2.25"
    "D42" = "121.89"
    "D43" = "21.15"
    "D44" = "0.0297"
    "D45" = "1.964.04"
    "D46" = "2.98"
    "D47" = "2.10"
    "D50" = "5.26"
    "D51" = "79.40"
}

$otherUpdates = @{
    "E2" = "  +4.35%  "
    "E3" = "  +2.46%  "
    "E4" = "  +0.16%  "
    "E5" = "  +1.14%  "
    "E6" = "  +2.38%  "
    "E7" = "  +1.53%  "
    "E8" = "  +0.06%  "
    "E9" = "  +2.46%  "
    "E10" = "  +6.80%  "
    "E11" = "  +0.93%  "
    "E12" = "  +1.07%  "
    "E13" = "  +0.88%  "
    "E14" = "  +1.16%  "
    "E15" = "  +2.66%  "
    "E16" = "  +3.14%  "
    "E17" = "  +0.42%  "
    "E18" = "  +4.45%  "
    "E19" = "  +4.28%  "
    "E20" = "  +3.24%  "
    "E21" = "  +1.32%  "
    "E22" = "  +2.36%  "
    "E23" = "  +6.45%  "
    "E24" = "  +2.61%  "
    "E25" = "  +3.24%  "
    "E26" = "  +1.35%  "
    "E27" = "  -0.06%  "
    "B28" = "Toncoin"
    "C28" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "E28" = "  -0.40%  "
    "B29" = "Cosmos"
    "C29" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "E29" = "  +4.14%  "
    "E30" = "  +6.73%  "
    "E31" = "  +8.30%  "
    "E32" = "  +0.05%  "
    "E33" = "  -3.30%  "
    "E34" = "  +3.19%  "
    "E35" = "  +2.50%  "
    "E36" = "  +0.29%  "
    "E37" = "  +5.31%  "
    "E38" = "  +3.46%  "
    "E39" = "  +3.40%  "
    "E40" = "  +1.85%  "
    "E41" = "  +1.86%  "
    "E42" = "  -3.04%  "
    "E43" = "  +3.34%  "
    "E44" = "  +2.41%  "
    "E45" = "  +1.37%  "
    "E46" = "  +1.94%  "
    "E47" = "  -0.54%  "
    "E48" = "  +0.66%  "
    "E49" = "  -0.18%  "
    "E50" = "  +9.37%  "
    "E51" = "  +3.30%  "
}

foreach ($cellRef in $priceUpdates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $priceUpdates[$cellRef]
}

foreach ($cellRef in $otherUpdates.Keys) {
    $ws.Range($cellRef).Value = $otherUpdates[$cellRef]
}

Write-Output ("Updated " + ($priceUpdates.Count + $otherUpdates.Count) + " cells")
